$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Fn1"
$ws.Cells.Item(2, 3).Value = "Itgb8"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 40.75339133333333
$ws.Cells.Item(2, 8).Value = 122.260174
$ws.Cells.Item(2, 9).Value = 0.02126536631186857
$ws.Cells.Item(2, 10).Value = 0.02126536631186857
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.027767
$ws.Cells.Item(2, 14).Value = 0.083301
$ws.Cells.Item(2, 15).Value = 0.002463719941166009
$ws.Cells.Item(2, 16).Value = 0.002463719941166009
$ws.Cells.Item(2, 17).Value = 1.131599417152667
$ws.Cells.Item(2, 18).Value = 10.184394754374
$ws.Cells.Item(2, 19).Value = 0.00005239190703875046
$ws.Cells.Item(2, 20).Value = 0.00005239190703875046

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Fn1"
$ws.Cells.Item(3, 3).Value = "Itgb8"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 40.75339133333333
$ws.Cells.Item(3, 8).Value = 122.260174
$ws.Cells.Item(3, 9).Value = 0.02126536631186857
$ws.Cells.Item(3, 10).Value = 0.02126536631186857
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 3.987076
$ws.Cells.Item(3, 14).Value = 11.961228
$ws.Cells.Item(3, 15).Value = 0.3537666527944829
$ws.Cells.Item(3, 16).Value = 0.3537666527944829
$ws.Cells.Item(3, 17).Value = 162.4868685037413
$ws.Cells.Item(3, 18).Value = 1462.381816533672
$ws.Cells.Item(3, 19).Value = 0.007522977460598301
$ws.Cells.Item(3, 20).Value = 0.007522977460598301

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Fn1"
$ws.Cells.Item(4, 3).Value = "Itgb8"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 40.75339133333333
$ws.Cells.Item(4, 8).Value = 122.260174
$ws.Cells.Item(4, 9).Value = 0.02126536631186857
$ws.Cells.Item(4, 10).Value = 0.02126536631186857
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 7.255512666666667
$ws.Cells.Item(4, 14).Value = 21.766538
$ws.Cells.Item(4, 15).Value = 0.643769627264351
$ws.Cells.Item(4, 16).Value = 0.643769627264351
$ws.Cells.Item(4, 17).Value = 295.6867470286236
$ws.Cells.Item(4, 18).Value = 2661.180723257612
$ws.Cells.Item(4, 19).Value = 0.01368999694423152
$ws.Cells.Item(4, 20).Value = 0.01368999694423152

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Fn1"
$ws.Cells.Item(5, 3).Value = "Itgb8"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 1689.289306666667
$ws.Cells.Item(5, 8).Value = 5067.86792
$ws.Cells.Item(5, 9).Value = 0.8814813868902838
$ws.Cells.Item(5, 10).Value = 0.8814813868902838
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.027767
$ws.Cells.Item(5, 14).Value = 0.083301
$ws.Cells.Item(5, 15).Value = 0.002463719941166009
$ws.Cells.Item(5, 16).Value = 0.002463719941166009
$ws.Cells.Item(5, 17).Value = 46.90649617821333
$ws.Cells.Item(5, 18).Value = 422.15846560392
$ws.Cells.Item(5, 19).Value = 0.002171723270648262
$ws.Cells.Item(5, 20).Value = 0.002171723270648262

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Fn1"
$ws.Cells.Item(6, 3).Value = "Itgb8"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 1689.289306666667
$ws.Cells.Item(6, 8).Value = 5067.86792
$ws.Cells.Item(6, 9).Value = 0.8814813868902838
$ws.Cells.Item(6, 10).Value = 0.8814813868902838
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 3.987076
$ws.Cells.Item(6, 14).Value = 11.961228
$ws.Cells.Item(6, 15).Value = 0.3537666527944829
$ws.Cells.Item(6, 16).Value = 0.3537666527944829
$ws.Cells.Item(6, 17).Value = 6735.324851667307
$ws.Cells.Item(6, 18).Value = 60617.92366500576
$ws.Cells.Item(6, 19).Value = 0.3118387197408143
$ws.Cells.Item(6, 20).Value = 0.3118387197408143

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Fn1"
$ws.Cells.Item(7, 3).Value = "Itgb8"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 1689.289306666667
$ws.Cells.Item(7, 8).Value = 5067.86792
$ws.Cells.Item(7, 9).Value = 0.8814813868902838
$ws.Cells.Item(7, 10).Value = 0.8814813868902838
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 7.255512666666667
$ws.Cells.Item(7, 14).Value = 21.766538
$ws.Cells.Item(7, 15).Value = 0.643769627264351
$ws.Cells.Item(7, 16).Value = 0.643769627264351
$ws.Cells.Item(7, 17).Value = 12256.65996218455
$ws.Cells.Item(7, 18).Value = 110309.939659661
$ws.Cells.Item(7, 19).Value = 0.5674709438788212
$ws.Cells.Item(7, 20).Value = 0.5674709438788212

# Row 8
$ws.Cells.Item(8, 1).Value = "M2"
$ws.Cells.Item(8, 2).Value = "Fn1"
$ws.Cells.Item(8, 3).Value = "Itgb8"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 100.9654023333333
$ws.Cells.Item(8, 8).Value = 302.896207
$ws.Cells.Item(8, 9).Value = 0.05268435816499466
$ws.Cells.Item(8, 10).Value = 0.05268435816499466
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.027767
$ws.Cells.Item(8, 14).Value = 0.083301
$ws.Cells.Item(8, 15).Value = 0.002463719941166009
$ws.Cells.Item(8, 16).Value = 0.002463719941166009
$ws.Cells.Item(8, 17).Value = 2.803506326589666
$ws.Cells.Item(8, 18).Value = 25.231556939307
$ws.Cells.Item(8, 19).Value = 0.0001297995037986296
$ws.Cells.Item(8, 20).Value = 0.0001297995037986296

# Row 9
$ws.Cells.Item(9, 1).Value = "M2"
$ws.Cells.Item(9, 2).Value = "Fn1"
$ws.Cells.Item(9, 3).Value = "Itgb8"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 100.9654023333333
$ws.Cells.Item(9, 8).Value = 302.896207
$ws.Cells.Item(9, 9).Value = 0.05268435816499466
$ws.Cells.Item(9, 10).Value = 0.05268435816499466
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 3.987076
$ws.Cells.Item(9, 14).Value = 11.961228
$ws.Cells.Item(9, 15).Value = 0.3537666527944829
$ws.Cells.Item(9, 16).Value = 0.3537666527944829
$ws.Cells.Item(9, 17).Value = 402.5567324735773
$ws.Cells.Item(9, 18).Value = 3623.010592262196
$ws.Cells.Item(9, 19).Value = 0.01863796904265585
$ws.Cells.Item(9, 20).Value = 0.01863796904265585

# Row 10
$ws.Cells.Item(10, 1).Value = "M2"
$ws.Cells.Item(10, 2).Value = "Fn1"
$ws.Cells.Item(10, 3).Value = "Itgb8"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 100.9654023333333
$ws.Cells.Item(10, 8).Value = 302.896207
$ws.Cells.Item(10, 9).Value = 0.05268435816499466
$ws.Cells.Item(10, 10).Value = 0.05268435816499466
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 7.255512666666667
$ws.Cells.Item(10, 14).Value = 21.766538
$ws.Cells.Item(10, 15).Value = 0.643769627264351
$ws.Cells.Item(10, 16).Value = 0.643769627264351
$ws.Cells.Item(10, 17).Value = 732.5557555245963
$ws.Cells.Item(10, 18).Value = 6593.001799721366
$ws.Cells.Item(10, 19).Value = 0.03391658961854018
$ws.Cells.Item(10, 20).Value = 0.03391658961854018

# Row 11
$ws.Cells.Item(11, 1).Value = "sCs"
$ws.Cells.Item(11, 2).Value = "Fn1"
$ws.Cells.Item(11, 3).Value = "Itgb8"
$ws.Cells.Item(11, 4).Value = "ECs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 85.41274733333334
$ws.Cells.Item(11, 8).Value = 256.238242
$ws.Cells.Item(11, 9).Value = 0.04456888863285297
$ws.Cells.Item(11, 10).Value = 0.04456888863285297
$ws.Cells.Item(11, 11).Value = 1
$ws.Cells.Item(11, 12).Value = 0.3333333333333333
$ws.Cells.Item(11, 13).Value = 0.027767
$ws.Cells.Item(11, 14).Value = 0.083301
$ws.Cells.Item(11, 15).Value = 0.002463719941166009
$ws.Cells.Item(11, 16).Value = 0.002463719941166009
$ws.Cells.Item(11, 17).Value = 2.371655755204667
$ws.Cells.Item(11, 18).Value = 21.344901796842
$ws.Cells.Item(11, 19).Value = 0.0001098052596803669
$ws.Cells.Item(11, 20).Value = 0.0001098052596803669

# Row 12
$ws.Cells.Item(12, 1).Value = "sCs"
$ws.Cells.Item(12, 2).Value = "Fn1"
$ws.Cells.Item(12, 3).Value = "Itgb8"
$ws.Cells.Item(12, 4).Value = "FAPs"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 85.41274733333334
$ws.Cells.Item(12, 8).Value = 256.238242
$ws.Cells.Item(12, 9).Value = 0.04456888863285297
$ws.Cells.Item(12, 10).Value = 0.04456888863285297
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 3.987076
$ws.Cells.Item(12, 14).Value = 11.961228
$ws.Cells.Item(12, 15).Value = 0.3537666527944829
$ws.Cells.Item(12, 16).Value = 0.3537666527944829
$ws.Cells.Item(12, 17).Value = 340.5471149867974
$ws.Cells.Item(12, 18).Value = 3064.924034881176
$ws.Cells.Item(12, 19).Value = 0.01576698655041447
$ws.Cells.Item(12, 20).Value = 0.01576698655041447

# Row 13
$ws.Cells.Item(13, 1).Value = "sCs"
$ws.Cells.Item(13, 2).Value = "Fn1"
$ws.Cells.Item(13, 3).Value = "Itgb8"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 85.41274733333334
$ws.Cells.Item(13, 8).Value = 256.238242
$ws.Cells.Item(13, 9).Value = 0.04456888863285297
$ws.Cells.Item(13, 10).Value = 0.04456888863285297
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 7.255512666666667
$ws.Cells.Item(13, 14).Value = 21.766538
$ws.Cells.Item(13, 15).Value = 0.643769627264351
$ws.Cells.Item(13, 16).Value = 0.643769627264351
$ws.Cells.Item(13, 17).Value = 619.7132701717996
$ws.Cells.Item(13, 18).Value = 5577.419431546196
$ws.Cells.Item(13, 19).Value = 0.02869209682275813
$ws.Cells.Item(13, 20).Value = 0.02869209682275813
